# Insert a new data row at row 145 (pushes the existing rows 145-225 down
# to 146-226, matching the rest of the A:T columns which are identical for
# every record in this sheet), then populate the new row with the values
# for the added "Fruta / hortaliza, semanal" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(145).Insert()

$ws.Range("A145").Value = 9
$ws.Range("B145").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C145").Value = "Metropolitana"
$ws.Range("D145").Value = 44875
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100101
$ws.Range("H145").Value = "Berries"
$ws.Range("I145").Value = 100101001
$ws.Range("J145").Value = "Arándano (blue)"
$ws.Range("K145").Value = "Sin especificar"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 510
$ws.Range("N145").Value = 6500
$ws.Range("O145").Value = 7000
$ws.Range("P145").Value = 6775
$ws.Range("Q145").Value = '$/bandeja 2 kilos'
$ws.Range("R145").Value = "Provincia de Linares"
$ws.Range("S145").Value = 3388
$ws.Range("T145").Value = 2
